# Edit script: update the "Förändrad" (changed) date column and append a new
# cleaning-notice row to the "Avverkningsanmälningar" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Column C ("Förändrad") holds the same date (45190) for every data row
#    (rows 2-533). The workbook was regenerated and that date moved on two
#    days to 45192 for all of them.
# ---------------------------------------------------------------------------
$lastRow = 533
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45192
}

# ---------------------------------------------------------------------------
# 2. Row 533 gains an explicit row height (matching every other data row).
# ---------------------------------------------------------------------------
$ws.Rows.Item(533).RowHeight = 15

# ---------------------------------------------------------------------------
# 3. Append a brand-new record as row 534.
# ---------------------------------------------------------------------------
$newRow = 534

$ws.Cells.Item($newRow, 1).Value = "A 44983-2023"

$ws.Cells.Item($newRow, 2).Value = 45190
$ws.Cells.Item($newRow, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($newRow, 3).Value = 45192
$ws.Cells.Item($newRow, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($newRow, 4).Value = "NORRBOTTENS LÄN"
$ws.Cells.Item($newRow, 5).Value = "GÄLLIVARE"
$ws.Cells.Item($newRow, 6).Value = "SCA"

$ws.Cells.Item($newRow, 7).Value = 2.6
$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 0
$ws.Cells.Item($newRow, 11).Value = 0
$ws.Cells.Item($newRow, 12).Value = 0
$ws.Cells.Item($newRow, 13).Value = 0
$ws.Cells.Item($newRow, 14).Value = 0
$ws.Cells.Item($newRow, 15).Value = 0
$ws.Cells.Item($newRow, 16).Value = 0
$ws.Cells.Item($newRow, 17).Value = 0

# Column R (Artnamn) keeps the wrap-text style used throughout the sheet,
# even though this new record has no species listed.
$ws.Cells.Item($newRow, 18).WrapText = $true
